# Generate Report for Handback
#
# This script reproduces the "handback" report-generation edit:
#  - Overview/zh-cn/de-de "Status" cells flip from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - per-language "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns get populated (with a hyperlink on
#    the target-file cell, matching the existing source-file hyperlink)
#  - a handful of columns are widened to fit the newly-populated content

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdFile1 = "0f81beb4-be80-4c40-a3b8-2029c2bd042e.md"
$mdFile2 = "a2d9c3fc-f13d-4169-9db8-044cfd773835.md"

$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/911a26f5e2311115f28afe9d310ca293cc01614c/e2e/0f81beb4-be80-4c40-a3b8-2029c2bd042e.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/911a26f5e2311115f28afe9d310ca293cc01614c/e2e/a2d9c3fc-f13d-4169-9db8-044cfd773835.md"

$zhXlf1 = "0f81beb4-be80-4c40-a3b8-2029c2bd042e.2e39f6d35af0355e0c36115847234163a797ed0d.zh-cn.xlf"
$zhXlf2 = "a2d9c3fc-f13d-4169-9db8-044cfd773835.cd1f390e900ad256c676384f8ba67b364e1b7765.zh-cn.xlf"
$deXlf1 = "0f81beb4-be80-4c40-a3b8-2029c2bd042e.2e39f6d35af0355e0c36115847234163a797ed0d.de-de.xlf"
$deXlf2 = "a2d9c3fc-f13d-4169-9db8-044cfd773835.cd1f390e900ad256c676384f8ba67b364e1b7765.de-de.xlf"

$zhHandbackTime = "2016-08-18 10:33:41"
$deHandbackTime = "2016-08-18 10:33:49"

# blue hyperlink-like font color used elsewhere in the workbook (FF6495ED)
$hyperlinkColor = 15570276

# ------------------------------------------------------------------
# 1) Overview sheet: flip the per-language status cells.
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

foreach ($addr in @("E2", "F2", "E3", "F3")) {
    $cell = $wsOverview.Range($addr)
    if ($cell.Text -eq $statusOld) {
        $cell.Value = $statusNew
    }
}

$wsOverview.Columns.Item(5).ColumnWidth = 29.1667   # E
$wsOverview.Columns.Item(6).ColumnWidth = 29.1667   # F

# ------------------------------------------------------------------
# 2) Per-language detail sheets (zh-cn, de-de).
# ------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; Xlf1 = $zhXlf1; Xlf2 = $zhXlf2; HandbackTime = $zhHandbackTime },
    @{ Sheet = "de-de"; Xlf1 = $deXlf1; Xlf2 = $deXlf2; HandbackTime = $deHandbackTime }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status column (C) for both rows.
    foreach ($addr in @("C2", "C3")) {
        $cell = $ws.Range($addr)
        if ($cell.Text -eq $statusOld) {
            $cell.Value = $statusNew
        }
    }

    # Row 2 -> 0f81beb4... file; Row 3 -> a2d9c3fc... file.
    $rows = @(
        @{ Row = 2; Md = $mdFile1; Url = $url1; Xlf = $lang.Xlf1 },
        @{ Row = 3; Md = $mdFile2; Url = $url2; Xlf = $lang.Xlf2 }
    )

    foreach ($r in $rows) {
        $iCell = $ws.Cells.Item($r.Row, 9)   # I: Latest Target File
        $jCell = $ws.Cells.Item($r.Row, 10)  # J: Latest Handback File
        $kCell = $ws.Cells.Item($r.Row, 11)  # K: Latest Handback DateTime

        $ws.Hyperlinks.Add($iCell, $r.Url, "", "", $r.Md)
        $iCell.Font.Underline = 2
        $iCell.Font.Color = $hyperlinkColor

        $jCell.Value = $r.Xlf
        $kCell.Value = $lang.HandbackTime
    }

    # Widen columns C (Status), I (Latest Target File), J (Latest Handback File).
    $ws.Columns.Item(3).ColumnWidth = 29.1667
    $ws.Columns.Item(9).ColumnWidth = 39.1667
    $ws.Columns.Item(10).ColumnWidth = 39.1667
}
